$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "densee_com_v3" (Q) and "densee_com_v5" (R) columns
# for rows 8 through 27.
$ws.Range("Q8").Value = 51
$ws.Range("R8").Value = 33

$ws.Range("Q9").Value = 57
$ws.Range("R9").Value = 47

$ws.Range("Q10").Value = 62
$ws.Range("R10").Value = 53

$ws.Range("Q11").Value = 67
$ws.Range("R11").Value = 53

$ws.Range("Q12").Value = 70
$ws.Range("R12").Value = 61

$ws.Range("Q13").Value = 71
$ws.Range("R13").Value = 57

$ws.Range("Q14").Value = 73
$ws.Range("R14").Value = 62

$ws.Range("Q15").Value = 71
$ws.Range("R15").Value = 58

$ws.Range("Q16").Value = 72
$ws.Range("R16").Value = 62

$ws.Range("Q17").Value = 74
$ws.Range("R17").Value = 62

$ws.Range("Q18").Value = 67
$ws.Range("R18").Value = 65

$ws.Range("Q19").Value = 74
$ws.Range("R19").Value = 64

$ws.Range("Q20").Value = 75
$ws.Range("R20").Value = 65

$ws.Range("Q21").Value = 76
$ws.Range("R21").Value = 62

$ws.Range("Q22").Value = 76
$ws.Range("R22").Value = 64

$ws.Range("Q23").Value = 77
$ws.Range("R23").Value = 62

$ws.Range("Q24").Value = 76
$ws.Range("R24").Value = 64

$ws.Range("Q25").Value = 76
$ws.Range("R25").Value = 65

$ws.Range("Q26").Value = 76
$ws.Range("R26").Value = 64

$ws.Range("Q27").Value = 76
$ws.Range("R27").Value = 64

# Scroll the view back so A1 is the top-left cell again and leave the
# active selection on R1 (matches the author's final on-screen state).
$ws.Range("A1").Select() | Out-Null
$ws.Range("R1").Select() | Out-Null
